$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 149.5
$ws.Range("I6").Value = 149.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 448.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -336.5
$ws.Range("N6").ClearContents()
$ws.Range("H15").Value = 929.48
$ws.Range("I15").Value = 929.48
$ws.Range("K15").Value = 2788.44
$ws.Range("M15").Value = -2619.44
$ws.Range("H21").Value = 99999
$ws.Range("J21").Value = 99999
$ws.Range("L21").Value = 99999
$ws.Range("N21").Value = -100935
$ws.Range("H23").Value = 99999
$ws.Range("J23").Value = 99999
$ws.Range("L23").Value = 99999
$ws.Range("N23").Value = -100467
$ws.Range("H40").Value = 4668.1577
$ws.Range("I40").Value = 3789.4
$ws.Range("J40").Value = 4982
$ws.Range("K40").Value = 3789.4
$ws.Range("L40").Value = 4982
$ws.Range("M40").Value = -3614.4
$ws.Range("N40").Value = -5332
$ws.Range("H53").Value = 2020.9565
$ws.Range("J53").Value = 2720
$ws.Range("L53").Value = 2720
$ws.Range("N53").Value = -3994
$ws.Range("H62").Value = 100027.73
$ws.Range("I62").Value = 133346.75
$ws.Range("K62").Value = 133346.75
$ws.Range("M62").Value = -132722.75
$ws.Range("H65").Value = 100027.73
$ws.Range("I65").Value = 133346.75
$ws.Range("K65").Value = 666733.75
$ws.Range("M65").Value = -663613.75
$ws.Range("H86").Value = 2746.1904
$ws.Range("I86").Value = 2191.4285
$ws.Range("J86").Value = 3855.7144
$ws.Range("K86").Value = 2191.4285
$ws.Range("L86").Value = 3855.7144
$ws.Range("M86").Value = -1068.4285
$ws.Range("N86").Value = -6101.7144
$ws.Range("H89").Value = 2746.1904
$ws.Range("I89").Value = 2191.4285
$ws.Range("J89").Value = 3855.7144
$ws.Range("K89").Value = 10957.1425
$ws.Range("L89").Value = 19278.572
$ws.Range("M89").Value = -5341.1425
$ws.Range("N89").Value = -30510.572
$ws.Range("H98").Value = 2107.7273
$ws.Range("I98").Value = 1641.8334
$ws.Range("J98").Value = 4204.25
$ws.Range("K98").Value = 1641.8334
$ws.Range("L98").Value = 4204.25
$ws.Range("M98").Value = -143.8334
$ws.Range("N98").Value = -7200.25
$ws.Range("H112").Value = 219585.02
$ws.Range("I112").Value = 3899.5
$ws.Range("J112").Value = 240126.5
$ws.Range("K112").Value = 11698.5
$ws.Range("L112").Value = 720379.5
$ws.Range("M112").Value = -10590.5
$ws.Range("N112").Value = -722595.5
$ws.Range("H116").Value = 4169.9165
$ws.Range("I116").Value = 3707.4
$ws.Range("K116").Value = 3707.4
$ws.Range("M116").Value = -265.4000000000001
$ws.Range("H122").Value = 2107.7273
$ws.Range("I122").Value = 1641.8334
$ws.Range("J122").Value = 4204.25
$ws.Range("K122").Value = 4925.5002
$ws.Range("L122").Value = 12612.75
$ws.Range("M122").Value = -2475.5002
$ws.Range("N122").Value = -17512.75
$ws.Range("H131").Value = 12300.6
$ws.Range("I131").Value = 14438.25
$ws.Range("K131").Value = 43314.75
$ws.Range("M131").Value = -38274.75
$ws.Range("H137").Value = 5868.1514
$ws.Range("J137").Value = 8714.143
$ws.Range("L137").Value = 26142.429
$ws.Range("N137").Value = -31242.429
$ws.Range("H141").Value = 6024.143
$ws.Range("I141").Value = 4436
$ws.Range("K141").Value = 13308
$ws.Range("M141").Value = -8128
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4054.9268
$ws.Range("I61").Value = 3993.8125
$ws.Range("K61").Value = 3993.8125
$ws.Range("M61").Value = -3781.8125
$ws.Range("H63").Value = 4097.5713
$ws.Range("I63").Value = 3779.5
$ws.Range("K63").Value = 3779.5
$ws.Range("M63").Value = -3093.5
$ws.Range("H66").Value = 4097.5713
$ws.Range("I66").Value = 3779.5
$ws.Range("K66").Value = 18897.5
$ws.Range("M66").Value = -15465.5
$ws.Range("H102").Value = 1880.4166
$ws.Range("I102").Value = 1456.5
$ws.Range("K102").Value = 1456.5
$ws.Range("M102").Value = 165.5
$ws.Range("H122").Value = 2297.0833
$ws.Range("I122").Value = 2219.3
$ws.Range("K122").Value = 6657.900000000001
$ws.Range("M122").Value = -4207.900000000001
$ws.Range("H136").Value = 4054.9268
$ws.Range("I136").Value = 3993.8125
$ws.Range("K136").Value = 11981.4375
$ws.Range("M136").Value = -9431.4375
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 766
$ws.Range("I22").Value = 999.5
$ws.Range("K22").Value = 999.5
$ws.Range("M22").Value = -826.5
$ws.Range("H24").Value = 1374.25
$ws.Range("I24").Value = 1374.25
$ws.Range("K24").Value = 1374.25
$ws.Range("M24").Value = -1139.25
$ws.Range("H94").Value = 1643.8928
$ws.Range("I94").Value = 1004.4091
$ws.Range("K94").Value = 1004.4091
$ws.Range("M94").Value = -553.4091
$ws.Range("H134").Value = 5606.5386
$ws.Range("I134").Value = 5740.4165
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 17221.2495
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -14686.2495
$ws.Range("N134").Value = -17070
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 505000000
$ws.Range("J4").Value = 505000000
$ws.Range("L4").Value = 505000000
$ws.Range("N4").Value = -505000224
$ws.Range("H7").Value = 239
$ws.Range("I7").Value = 282.25
$ws.Range("J7").Value = 152.5
$ws.Range("K7").Value = 282.25
$ws.Range("L7").Value = 152.5
$ws.Range("M7").Value = -169.25
$ws.Range("N7").Value = -378.5
$ws.Range("H31").Value = 1970.2115
$ws.Range("I31").Value = 1399.6522
$ws.Range("K31").Value = 1399.6522
$ws.Range("M31").Value = -1104.6522
$ws.Range("H34").Value = 1970.2115
$ws.Range("I34").Value = 1399.6522
$ws.Range("K34").Value = 1399.6522
$ws.Range("M34").Value = -1197.6522
$ws.Range("H58").Value = 3998.2
$ws.Range("I58").Value = 3653.2856
$ws.Range("K58").Value = 3653.2856
$ws.Range("M58").Value = -3450.2856
$ws.Range("H92").Value = 65581.664
$ws.Range("J92").Value = 65581.664
$ws.Range("L92").Value = 65581.664
$ws.Range("N92").Value = -70573.664
$ws.Range("H99").Value = 2951.3333
$ws.Range("I99").Value = 2784.7144
$ws.Range("K99").Value = 2784.7144
$ws.Range("M99").Value = -1286.7144
$ws.Range("H107").Value = 327.35898
$ws.Range("I107").Value = 189.04167
$ws.Range("K107").Value = 189.04167
$ws.Range("M107").Value = 1730.95833
$ws.Range("H126").Value = 2951.3333
$ws.Range("I126").Value = 2784.7144
$ws.Range("K126").Value = 8354.143199999999
$ws.Range("M126").Value = -5884.143199999999
$ws.Range("H136").Value = 3998.2
$ws.Range("I136").Value = 3653.2856
$ws.Range("K136").Value = 10959.8568
$ws.Range("M136").Value = -8409.856800000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3897.6858
$ws.Range("I131").Value = 10924.25
$ws.Range("K131").Value = 32772.75
$ws.Range("M131").Value = -27732.75
$ws.Range("H132").Value = 1534.3
$ws.Range("I132").Value = 909.75
$ws.Range("K132").Value = 8187.75
$ws.Range("M132").Value = -5657.75
$ws.Range("H140").Value = 2371.3635
$ws.Range("I140").Value = 802.75
$ws.Range("J140").Value = 4784.615
$ws.Range("K140").Value = 2408.25
$ws.Range("L140").Value = 14353.845
$ws.Range("M140").Value = 2771.75
$ws.Range("N140").Value = -24713.845
$ws.Range("H141").Value = 12828.533
$ws.Range("I141").Value = 11648.308
$ws.Range("K141").Value = 34944.924
$ws.Range("M141").Value = -29764.924
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 43592
$ws.Range("J5").Value = 9000
$ws.Range("L5").Value = 9000
$ws.Range("N5").Value = -9224
$ws.Range("H102").Value = 64633.625
$ws.Range("I102").Value = 2472.3076
$ws.Range("J102").Value = 333999.34
$ws.Range("K102").Value = 2472.3076
$ws.Range("L102").Value = 333999.34
$ws.Range("M102").Value = -850.3076000000001
$ws.Range("N102").Value = -337243.34
$ws.Range("H122").Value = 1183.44
$ws.Range("I122").Value = 1246.4445
$ws.Range("J122").Value = 1021.4286
$ws.Range("K122").Value = 3739.3335
$ws.Range("L122").Value = 3064.2858
$ws.Range("M122").Value = -1289.3335
$ws.Range("N122").Value = -7964.2858
$ws.Range("H126").Value = 8973.6
$ws.Range("I126").Value = 5687
$ws.Range("J126").Value = 11164.667
$ws.Range("K126").Value = 17061
$ws.Range("L126").Value = 33494.001
$ws.Range("M126").Value = -14591
$ws.Range("N126").Value = -38434.001
$ws.Range("H132").Value = 1622.3889
$ws.Range("I132").Value = 1647.9375
$ws.Range("K132").Value = 4943.8125
$ws.Range("M132").Value = -2413.8125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 19776.4
$ws.Range("I122").Value = 19442.5
$ws.Range("J122").Value = 19999
$ws.Range("K122").Value = 58327.5
$ws.Range("L122").Value = 59997
$ws.Range("M122").Value = -55877.5
$ws.Range("N122").Value = -64897
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 3710.5715
$ws.Range("J13").Value = 2999.5
$ws.Range("L13").Value = 2999.5
$ws.Range("N13").Value = -3279.5
$ws.Range("H136").Value = 1752.0625
$ws.Range("I136").Value = 1647
$ws.Range("K136").Value = 4941
$ws.Range("M136").Value = -2391
$ws.Range("H138").Value = 99573
$ws.Range("J138").Value = 99573
$ws.Range("L138").Value = 99573
$ws.Range("N138").Value = -109853
Write-Output "done"
